$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-03-10 Sunday"; New = "2024-03-11 Monday" },
    @{ Old = "20÷6="; New = "77÷4=" },
    @{ Old = "71÷7="; New = "72÷2=" },
    @{ Old = "52÷4="; New = "71÷6=" },
    @{ Old = "22÷3="; New = "44÷8=" },
    @{ Old = "78÷4="; New = "84÷2=" },
    @{ Old = "57÷6="; New = "40÷9=" },
    @{ Old = "31÷4="; New = "73÷7=" },
    @{ Old = "45÷3="; New = "65÷3=" },
    @{ Old = "55÷2="; New = "21÷9=" },
    @{ Old = "48÷8="; New = "16÷8=" },
    @{ Old = "74÷6="; New = "14÷7=" },
    @{ Old = "98÷9="; New = "33÷6=" },
    @{ Old = "51÷4="; New = "42÷7=" },
    @{ Old = "81÷7="; New = "17÷2=" },
    @{ Old = "32÷4="; New = "21÷9=" },
    @{ Old = "46÷9="; New = "14÷6=" },
    @{ Old = "86÷4="; New = "19÷8=" },
    @{ Old = "38÷6="; New = "92÷8=" },
    @{ Old = "68÷8="; New = "54÷5=" },
    @{ Old = "84÷5="; New = "27÷8=" },
    @{ Old = "73÷3="; New = "28÷7=" },
    @{ Old = "25÷6="; New = "34÷7=" },
    @{ Old = "63÷3="; New = "43÷8=" },
    @{ Old = "46÷2="; New = "35÷8=" },
    @{ Old = "24÷3="; New = "35÷5=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
